$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 103
$ws.Range("J2").Value = 399
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 128
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 64
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 7
$ws.Range("T2").Value = 72
$ws.Range("U2").Value = 6
$ws.Range("V2").Value = 618
$ws.Range("X2").Value = 639
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 12
$ws.Range("AA2").Value = 6
